$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 29470
$ws.Range("C3").Value = 3953
$ws.Range("D3").Value = 5138
$ws.Range("B4").Value = 13086
$ws.Range("C4").Value = 1570
$ws.Range("D4").Value = 1792
$ws.Range("B5").Value = 46082
$ws.Range("C5").Value = 3362
$ws.Range("D5").Value = 4258
$ws.Range("B6").Value = 738
$ws.Range("C6").Value = 317
$ws.Range("D6").Value = 29
$ws.Range("B7").Value = 27107
$ws.Range("C7").Value = 5186
$ws.Range("D7").Value = 4547
$ws.Range("B8").Value = 3279
$ws.Range("C8").Value = 754
$ws.Range("D8").Value = 664
$ws.Range("B9").Value = 3946
$ws.Range("C9").Value = 647
$ws.Range("D9").Value = 572
$ws.Range("B10").Value = 1184
$ws.Range("C10").Value = 199
$ws.Range("D10").Value = 50
$ws.Range("B11").Value = 77
$ws.Range("C11").Value = 153
$ws.Range("D11").Value = 1
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("B13").Value = 720
$ws.Range("C13").Value = 205
$ws.Range("D13").Value = 173
$ws.Range("B14").Value = 1803
$ws.Range("C14").Value = 1032
$ws.Range("D14").Value = 732
$ws.Range("B15").Value = 3605
$ws.Range("C15").Value = 1414
$ws.Range("D15").Value = 676
$ws.Range("B16").Value = 2359
$ws.Range("C16").Value = 1032
$ws.Range("D16").Value = 478
$ws.Range("B17").Value = 1126
$ws.Range("C17").Value = 445
$ws.Range("D17").Value = 90
$ws.Range("B18").Value = 11099
$ws.Range("C18").Value = 1735
$ws.Range("D18").Value = 2067
$ws.Range("B19").Value = 2126
$ws.Range("C19").Value = 531
$ws.Range("D19").Value = 411
$ws.Range("B20").Value = 13978
$ws.Range("C20").Value = 1493
$ws.Range("D20").Value = 2469
$ws.Range("B21").Value = 247
$ws.Range("C21").Value = 270
$ws.Range("D21").Value = 11
$ws.Range("B22").Value = 12300
$ws.Range("C22").Value = 1433
$ws.Range("D22").Value = 1968
$ws.Range("B23").Value = 709
$ws.Range("C23").Value = 345
$ws.Range("D23").Value = 102
$ws.Range("B24").Value = 11120
$ws.Range("C24").Value = 2220
$ws.Range("D24").Value = 2077
$ws.Range("B25").Value = 45468
$ws.Range("C25").Value = 5531
$ws.Range("D25").Value = 5271
$ws.Range("B26").Value = 3357
$ws.Range("C26").Value = 1018
$ws.Range("D26").Value = 464
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("B28").Value = 3734
$ws.Range("C28").Value = 741
$ws.Range("D28").Value = 851
$ws.Range("B29").Value = 1815
$ws.Range("C29").Value = 219
$ws.Range("D29").Value = 385
$ws.Range("B30").Value = 9778
$ws.Range("C30").Value = 1894
$ws.Range("D30").Value = 1623
$ws.Range("B31").Value = 304
$ws.Range("C31").Value = 206
$ws.Range("D31").Value = 125
$ws.Range("B32").Value = 1206
$ws.Range("C32").Value = 1085
$ws.Range("D32").Value = 175
$ws.Range("B33").Value = 9784
$ws.Range("C33").Value = 2064
$ws.Range("D33").Value = 1917
$ws.Range("B34").Value = 7457
$ws.Range("C34").Value = 1808
$ws.Range("D34").Value = 1958
$ws.Range("B35").Value = 4049
$ws.Range("C35").Value = 516
$ws.Range("D35").Value = 954
$ws.Range("B36").Value = 32206
$ws.Range("C36").Value = 3701
$ws.Range("D36").Value = 3794
$ws.Range("B37").Value = 5547
$ws.Range("C37").Value = 1794
$ws.Range("D37").Value = 942
$ws.Range("B38").Value = 11609
$ws.Range("C38").Value = 1281
$ws.Range("D38").Value = 1960
$ws.Range("B39").Value = 438
$ws.Range("C39").Value = 490
$ws.Range("D39").Value = 102
$ws.Range("B40").Value = 1012
$ws.Range("C40").Value = 165
$ws.Range("D40").Value = 292
$ws.Range("B41").Value = 1826
$ws.Range("C41").Value = 234
$ws.Range("D41").Value = 130
$ws.Range("B42").Value = 7719
$ws.Range("C42").Value = 248
$ws.Range("D42").Value = 192
$ws.Range("B43").Value = 240
$ws.Range("C43").Value = 102
$ws.Range("D43").Value = 72
$ws.Range("B44").Value = 551
$ws.Range("C44").Value = 48
$ws.Range("D44").Value = 22
$ws.Range("B45").Value = 1743
$ws.Range("C45").Value = 159
$ws.Range("D45").Value = 73
$ws.Range("B46").Value = 2014
$ws.Range("C46").Value = 619
$ws.Range("D46").Value = 283
$ws.Range("B47").Value = 7839
$ws.Range("C47").Value = 2161
$ws.Range("D47").Value = 1611
$ws.Range("B48").Value = 20751
$ws.Range("C48").Value = 2077
$ws.Range("D48").Value = 3348
$ws.Range("B49").Value = 8911
$ws.Range("C49").Value = 2188
$ws.Range("D49").Value = 762
$ws.Range("B50").Value = 6919
$ws.Range("C50").Value = 731
$ws.Range("D50").Value = 1173
$ws.Range("B51").Value = 16641
$ws.Range("C51").Value = 2012
$ws.Range("D51").Value = 2093
$ws.Range("B52").Value = 2591
$ws.Range("C52").Value = 321
$ws.Range("D52").Value = 549
$ws.Range("B53").Value = 8519
$ws.Range("C53").Value = 1892
$ws.Range("D53").Value = 1685
$ws.Range("B54").Value = 1460
$ws.Range("C54").Value = 454
$ws.Range("D54").Value = 814
$ws.Range("B55").Value = 1256
$ws.Range("C55").Value = 905
$ws.Range("D55").Value = 177
$ws.Range("B56").Value = 2328
$ws.Range("C56").Value = 686
$ws.Range("D56").Value = 852
$ws.Range("B57").Value = 8487
$ws.Range("C57").Value = 3683
$ws.Range("D57").Value = 2059
$ws.Range("B58").Value = 12079
$ws.Range("C58").Value = 791
$ws.Range("D58").Value = 489
$ws.Range("B59").Value = 415233
$ws.Range("C59").Value = 65646
$ws.Range("D59").Value = 62954
